# "Generate Report for Handoff"
#
# Rows 5-10 of each sheet represent files that have just been handed off
# for localization (zh-cn / de-de). This updates:
#   - Overview sheet:  "Latest HO Xliff Generate Date" (col G) for rows 5-10
#   - zh-cn sheet:      "Priority" (col E) -> "ht" and
#                       "Latest Handoff Datetime" (col H) for rows 5-10
#   - de-de sheet:      "Priority" (col E) -> "ht" and
#                       "Latest Handoff Datetime" (col H) for rows 5-10

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

for ($row = 5; $row -le 10; $row++) {
    $wsOverview.Range("G$row").Value = "2016-11-03 20:35:57"

    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-11-03 20:35:44"

    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-11-03 20:35:57"
}
